$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)
$p.Range.ListFormat.RemoveNumbers()
